$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: ANGELICA MARIA GULFO BASTIDAS, period moved from 1708 -> 1707, Valor Mora 0 -> 5000000
$ws.Range("E16").Value = "1707"
$ws.Range("G16").Value = 5000000

# Row 17: was ANGELICA/1707 -> now ARLES MALDONADO WILCHES, 79908689, period 1707
$ws.Range("C17").Value = "79908689"
$ws.Range("D17").Value = "ARLES MALDONADO WILCHES"
$ws.Range("E17").Value = "1707"
$ws.Range("F17").Value = 42300
$ws.Range("G17").Value = 1057500

# Row 18: was ARLES/1708 -> now ANGELICA MARIA GULFO BASTIDAS, 52718112, period 1708
$ws.Range("C18").Value = "52718112"
$ws.Range("D18").Value = "ANGELICA MARIA GULFO BASTIDAS"
$ws.Range("E18").Value = "1708"
$ws.Range("F18").Value = 200000
$ws.Range("G18").Value = 5000000

# Row 19: ARLES MALDONADO WILCHES, period moved from 1707 -> 1708 (values unchanged)
$ws.Range("E19").Value = "1708"
